$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.509.34"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.741.97"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.00"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4451"
$ws.Range("E7").Value = "  +4.97%  "
$ws.Range("E8").Value = "  -2.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07397"
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.50"
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.077"
$ws.Range("E11").Value = "  -1.91%  "
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.43"
$ws.Range("E13").Value = "  -1.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.901"
$ws.Range("E14").Value = "  -2.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.075"
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("D16").Value = "1.740.89"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.36"
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("E18").Value = "  -1.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06375"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.81"
$ws.Range("E21").Value = "  -1.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.727"
$ws.Range("E22").Value = "  -2.87%  "
$ws.Range("D23").Value = "27.543.07"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.13"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.098"
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.85"
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.06"
$ws.Range("E27").Value = "  -0.89%  "
$ws.Range("D28").Value = "1.941.00"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "124.93"
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.034"
$ws.Range("E30").Value = "  -3.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.046"
$ws.Range("E31").Value = "  -5.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09070"
$ws.Range("E32").Value = "  +1.85%  "
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.365"
$ws.Range("E34").Value = "  -3.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02272"
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("E36").Value = "  -4.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06015"
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2060"
$ws.Range("E38").Value = "  -1.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6234"
$ws.Range("E39").Value = "  -1.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.885"
$ws.Range("E40").Value = "  -1.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.180"
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.372"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.707"
$ws.Range("E43").Value = "  -2.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.15"
$ws.Range("E44").Value = "  -1.40%  "
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5785"
$ws.Range("E46").Value = "  -1.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.77"
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.924"
$ws.Range("E48").Value = "  -2.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06838"
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.113"
$ws.Range("E50").Value = "  -4.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.29"
$ws.Range("E51").Value = "  -2.33%  "
